$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("기계")
$ws.Activate()

# Capture the existing comment text before the column shift moves its
# logical home (D2 "시작시간" -> C2 after the delete).
$existingComment = $ws.Range("D2").Comment
$commentText = $existingComment.Text()

$ws.Columns("C:C").Select()
$ws.Columns("C:C").Delete()

# The comment object stays anchored to the literal cell "D2", which after
# the shift is now the end-time cell. Re-home it onto C2 (start-time),
# matching where the original comment's data actually belongs.
$ws.Range("D2").Comment.Delete()
$ws.Range("C2").AddComment($commentText)

# Row 2 no longer contains the taller "맑은 고딕" runs that used to live in
# the deleted machine-name column, so Excel's row autofit shrinks it back
# down to the plain-Arial default height.
$ws.Rows("2:2").RowHeight = 12.75
